$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new diary entry row (row 15)
$ws.Range("A15").Value = "24 loka"
$ws.Range("B15").Value = "14.45-16.00, 17.15-18.30"
$ws.Range("C15").Value = "Firework puuhat jatkuivat"
$ws.Range("E15").Value = "Modernisoinnissa haastetta, ehkä pitäisi ensi kerralla enemmän kopioida esimerkkimoottorista asioita, ja sitten pikkuhiljaa muutella mutta nyt pitää mennä tällä. "
$ws.Range("D15").Value = "Ok, ihan hyvä meininki. Toisten koodin tutkimisestakin oppii paljon ja sitäkin pitää työelämää varten paljon tehdä."
$ws.Range("G15").Value = 2.5

# Match the existing look: B holds a time range (wrapped, time number format),
# C/D/E hold wrapped free text.
$ws.Range("B15").NumberFormat = "h:mm"
$ws.Range("B15").WrapText = $true
$ws.Range("C15").WrapText = $true
$ws.Range("D15").WrapText = $true
$ws.Range("E15").WrapText = $true

# Row grew to fit the wrapped text, same as the other multi-line rows above
$ws.Rows.Item(15).RowHeight = 72.5

# Update the selection to reflect the new active cell
$ws.Range("E15").Select()

$wb.Save()
